# The workbook tracks a pool of unused "names" (Sheet1) and a log of
# names that have already been used (the "used" sheet). This script
# takes the first unused name off the top of Sheet1 (shifting the rest
# of the list up by one row) and appends a new record for it to the
# bottom of the "used" sheet, recording which file it was used for and
# when.

$wb = $excel.ActiveWorkbook

$namesSheet = $wb.Worksheets.Item("Sheet1")
$usedSheet  = $wb.Worksheets.Item("used")

# Grab the id that is about to be consumed (top of the names list).
$usedId = $namesSheet.Range("A1").Value2

# Remove it from the pool - this shifts every remaining row up by one,
# so the list stays contiguous starting at row 1.
$namesSheet.Rows.Item(1).Delete()

# Find the next free row at the bottom of the "used" log and record the
# newly consumed id along with the file it was used for and the
# timestamp it was used at.
$nextRow = $usedSheet.Cells.Item($usedSheet.Rows.Count, 1).End(-4162).Row + 1

$usedSheet.Cells.Item($nextRow, 1).Value = $usedId
$usedSheet.Cells.Item($nextRow, 2).Value = "ChatGPT Image 2026年1月18日 06_29_52.png"
$usedSheet.Cells.Item($nextRow, 3).Value = "2026-01-18 06:32:41"
